$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '51.167.42'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.79%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.941.99'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -1.57%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '375.97'
$ws.Range('D5').Style = "Normal"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '102.56'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -3.36%  '
$ws.Range('E7').Value = '  -1.76%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  -2.47%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.80'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -1.94%  '
$ws.Range('E11').Value = '  -0.88%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0840'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -0.54%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.401.69'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -1.87%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '17.96'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -4.02%  '
$ws.Range('E15').Value = '  -2.24%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.944.55'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -1.79%  '
$ws.Range('E17').Value = '  +0.22%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '51.099.10'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.07%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.16'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -6.67%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.13'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -3.88%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.59'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -3.41%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.0₃0956'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.37%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '263.04'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.32%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '68.22'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -1.86%  '
$ws.Range('E25').Value = '  +2.33%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.93'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +9.37%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.13'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +9.27%  '
$ws.Range('E28').Value = '  +5.18%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.167'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -2.34%  '
$ws.Range('B30').Value = 'Dai'
$ws.Range('C30').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.00'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('E31').Value = '  -1.60%  '
$ws.Range('E32').Value = '  -0.52%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '34.18'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -1.54%  '
$ws.Range('B34').Value = 'VeChain'
$ws.Range('C34').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0457'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.86%  '
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '50.72'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -1.32%  '
$ws.Range('E36').Value = '  -3.83%  '
$ws.Range('E37').Value = '  -0.17%  '
$ws.Range('E38').Value = '  -4.31%  '
$ws.Range('E39').Value = '  -1.64%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '16.45'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -5.99%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.115'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -1.65%  '
$ws.Range('E42').Value = '  -3.86%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '121.78'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -1.60%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '21.10'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -5.51%  '
$ws.Range('E45').Value = '  -1.51%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.273'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.52%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.33'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.66%  '
$ws.Range('E48').Value = '  -1.08%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.000.58'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -2.71%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0346'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -1.50%  '
$ws.Range('E51').Value = '  -3.17%  '
